# Apply updated crypto price/volume figures to columns D (Price) and E (Volume(1h)).
# Values that would otherwise be auto-parsed by Excel as numbers (single-dot
# decimals such as "555.78") are entered with a leading apostrophe so they stay
# plain text, matching the original inline-string cells in the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '65.150.20'
$ws.Range("E2").Value = '  +0.90%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '3.375.38'
$ws.Range("E3").Value = '  +0.48%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  -0.03%  '

# Row 5: BNB
$ws.Range("D5").Value = '''555.78'
$ws.Range("E5").Value = '  +0.17%  '

# Row 6: Solana
$ws.Range("D6").Value = '''174.67'
$ws.Range("E6").Value = '  -0.49%  '

# Row 7: XRP
$ws.Range("E7").Value = '  +1.98%  '

# Row 8: LidoStakedEther
$ws.Range("D8").Value = '3.363.44'
$ws.Range("E8").Value = '  +0.37%  '

# Row 9: USDC
$ws.Range("E9").Value = '  -0.14%  '

# Row 10: Dogecoin
$ws.Range("D10").Value = '''0.172'
$ws.Range("E10").Value = '  +5.68%  '

# Row 11: Cardano
$ws.Range("D11").Value = '''0.636'
$ws.Range("E11").Value = '  +1.08%  '

# Row 12: Avalanche
$ws.Range("D12").Value = '''53.69'
$ws.Range("E12").Value = '  -1.36%  '

# Row 13: ShibaInu
$ws.Range("D13").Value = '''0.0000279'
$ws.Range("E13").Value = '  +1.88%  '

# Row 14: Polkadot
$ws.Range("D14").Value = '''9.17'
$ws.Range("E14").Value = '  +0.85%  '

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = '3.914.74'
$ws.Range("E15").Value = '  +0.50%  '

# Row 16: Chainlink
$ws.Range("D16").Value = '''18.33'
$ws.Range("E16").Value = '  -0.44%  '

# Row 17: TRON
$ws.Range("E17").Value = '  +0.85%  '

# Row 18: WrappedEther
$ws.Range("D18").Value = '3.363.85'
$ws.Range("E18").Value = '  +0.10%  '

# Row 19: WrappedBTC
$ws.Range("D19").Value = '65.018.93'
$ws.Range("E19").Value = '  +0.80%  '

# Row 20: Uniswap
$ws.Range("D20").Value = '''11.84'
$ws.Range("E20").Value = '  +0.08%  '

# Row 21: Polygon
$ws.Range("D21").Value = '''0.998'
$ws.Range("E21").Value = '  +1.24%  '

# Row 22: BitcoinCash
$ws.Range("E22").Value = '  -1.16%  '

# Row 23: Toncoin
$ws.Range("D23").Value = '''4.89'
$ws.Range("E23").Value = '  +2.19%  '

# Row 24: PancakeSwap
$ws.Range("D24").Value = '''4.08'
$ws.Range("E24").Value = '  -0.17%  '

# Row 25: InternetComputer(DFINITY)
$ws.Range("D25").Value = '''14.11'
$ws.Range("E25").Value = '  +5.86%  '

# Row 26: Litecoin
$ws.Range("D26").Value = '''87.79'
$ws.Range("E26").Value = '  +2.03%  '

# Row 27: ImmutableX
$ws.Range("E27").Value = '  +1.90%  '

# Row 28: RenderToken
$ws.Range("D28").Value = '''10.69'
$ws.Range("E28").Value = '  -1.99%  '

# Row 29: Filecoin
$ws.Range("D29").Value = '''8.71'
$ws.Range("E29").Value = '  -0.64%  '

# Row 30: EthereumClassic
$ws.Range("D30").Value = '''31.06'
$ws.Range("E30").Value = '  +3.12%  '

# Row 31: NEARProtocol
$ws.Range("D31").Value = '''6.54'
$ws.Range("E31").Value = '  -1.51%  '

# Row 32: OKB
$ws.Range("D32").Value = '''63.15'
$ws.Range("E32").Value = '  +7.13%  '

# Row 33: Cosmos
$ws.Range("D33").Value = '''11.45'
$ws.Range("E33").Value = '  -0.28%  '

# Row 34: Bittensor
$ws.Range("D34").Value = '''577.56'
$ws.Range("E34").Value = '  -0.88%  '

# Row 35: Hedera
$ws.Range("D35").Value = '''0.108'
$ws.Range("E35").Value = '  -0.51%  '

# Row 36: Dai
$ws.Range("E36").Value = '  -0.03%  '

# Row 37: Stacks
$ws.Range("D37").Value = '''3.63'
$ws.Range("E37").Value = '  +3.97%  '

# Row 38: Kaspa
$ws.Range("E38").Value = '  +1.54%  '

# Row 39: InjectiveProtocol
$ws.Range("D39").Value = '''35.71'
$ws.Range("E39").Value = '  +0.19%  '

# Row 40: TheGraph
$ws.Range("E40").Value = '  -0.87%  '

# Row 41: PEPE
$ws.Range("D41").Value = '0.0₃0739'
$ws.Range("E41").Value = '  -2.57%  '

# Row 42: Maker
$ws.Range("D42").Value = '3.089.79'
$ws.Range("E42").Value = '  -0.43%  '

# Row 43: VeChain
$ws.Range("D43").Value = '''0.0417'
$ws.Range("E43").Value = '  +1.26%  '

# Row 44: ThetaToken
$ws.Range("D44").Value = '''2.76'
$ws.Range("E44").Value = '  -1.43%  '

# Row 45: ApeXProtocol
$ws.Range("D45").Value = '''3.20'
$ws.Range("E45").Value = '  -0.54%  '

# Row 46: Fetch.AI
$ws.Range("D46").Value = '''2.45'
$ws.Range("E46").Value = '  -3.12%  '

# Row 47: Stellar
$ws.Range("D47").Value = '''0.134'
$ws.Range("E47").Value = '  +2.02%  '

# Row 48: FirstDigitalUSD
$ws.Range("D48").Value = '''0.999'
$ws.Range("E48").Value = '  +0.06%  '

# Row 49: Monero
$ws.Range("D49").Value = '''141.10'
$ws.Range("E49").Value = '  +4.24%  '

# Row 50: WEMIXToken
$ws.Range("E50").Value = '  -2.12%  '

# Row 51: THORChain
$ws.Range("D51").Value = '''8.30'
$ws.Range("E51").Value = '  -0.99%  '
